# Added PDQ drug info summary load test
#
# The existing "PDQPage" sheet (lung-prevention PDQ Cancer Info Summary
# load-test data) is duplicated. The duplicate becomes the new first sheet,
# "PDQCisPage", keeping the original lung-prevention rows. The original
# sheet object moves to the second position and is renamed "PDQDrugPage";
# its rows are replaced with new PDQ Drug Info Summary load-test data.
# "RightNav" remains the third/last sheet.

$wb = $excel.ActiveWorkbook

$orig = $wb.Worksheets.Item("PDQPage")

# Duplicate the sheet; the copy is placed directly before the original.
$orig.Copy($orig, $null)

$cisPage = $wb.Worksheets.Item(1)
$cisPage.Name = "PDQCisPage"

$drugPage = $wb.Worksheets.Item(2)
$drugPage.Name = "PDQDrugPage"

# Drop the old lung-prevention detail rows (rows 5-10), keep header + the
# first 3 data rows' shape, then overwrite with the new drug content.
$drugPage.Rows("5:10").Delete()

$drugPage.Range("A2").Value = "/about-cancer/treatment/drugs/acalabrutinib"
$drugPage.Range("A3").Value = "/about-cancer/treatment/drugs/recombinant-HPV-quadrivalent-vaccine"
$drugPage.Range("A4").Value = "/about-cancer/treatment/drugs/sorafenibtosylate"

$drugPage.Range("B2").Value = "PDQ Drug Info Summary"
$drugPage.Range("B3").Value = "PDQ Drug Info Summary"
$drugPage.Range("B4").Value = "PDQ Drug Info Summary"

# Restore selections on each sheet and make the new drug page the active tab.
$cisPage.Range("A11").Select()

$drugPage.Activate()
$drugPage.Range("A3").Select()
